$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Overview" (sheet1) - add row 3 for the new handoff entry
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$newMdName        = "373b367f-fe1c-4a86-83ed-56377b13f714oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdPath        = "e2e\373b367f-fe1c-4a86-83ed-56377b13f714oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdUrl         = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1101e1d1e53de022c93cc4a9abf0bbfc46be5bf6/e2e/373b367f-fe1c-4a86-83ed-56377b13f714oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$readyForHandoff  = "Ready for handoff"
$handoffDate1     = "2016-08-30 10:07:19"

# Expand the table (ListObject) by one row so ref/autoFilter grow to A1:G3
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Range("B3").Value = $newMdPath
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = $handoffDate1

# B3 mirrors B2's hyperlink styling (underline + cornflowerblue font)
$wsOverview.Range("B3").Font.Underline = 2
$wsOverview.Range("B3").Font.Color = 15570276

# G3 mirrors G2's date-time number format
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdUrl, "", "", $newMdPath) | Out-Null

# Column widths E/F grow on the Overview sheet too
$wsOverview.Columns.Item(5).ColumnWidth = 17
$wsOverview.Columns.Item(6).ColumnWidth = 17

# ------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - add row 3 for the new handoff entry
# ------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$newXlfZh = "373b367f-fe1c-4a86-83ed-56377b13f714ooooooooooooooooooooooooooooooooooooooo.582648cef18041beebe4dd37345fb26f3865783d.zh-cn.xlf"
$xlfDateZh = "2016-08-30 10:07:04"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newMdName
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $xlfDateZh
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("M3").Value = "True"
$wsZh.Range("O3").Value = "False"

# A3 mirrors A2's hyperlink styling (underline + cornflowerblue font)
$wsZh.Range("A3").Font.Underline = 2
$wsZh.Range("A3").Font.Color = 15570276

# H3/K3 mirror the date-time number format used on row 2
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newMdUrl, "", "", $newMdName) | Out-Null

# Column width C grows on the zh-cn sheet too
$wsZh.Columns.Item(3).ColumnWidth = 17

# ------------------------------------------------------------------
# Sheet "de-de" (sheet3) - add row 3 for the new handoff entry
# ------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$newXlfDe = "373b367f-fe1c-4a86-83ed-56377b13f714ooooooooooooooooooooooooooooooooooooooo.582648cef18041beebe4dd37345fb26f3865783d.de-de.xlf"
$xlfDateDe = $handoffDate1

$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newMdName
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $xlfDateDe
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("M3").Value = "True"
$wsDe.Range("O3").Value = "False"

# A3 mirrors A2's hyperlink styling (underline + cornflowerblue font)
$wsDe.Range("A3").Font.Underline = 2
$wsDe.Range("A3").Font.Color = 15570276

# H3/K3 mirror the date-time number format used on row 2
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newMdUrl, "", "", $newMdName) | Out-Null

# Column width C grows on the de-de sheet too
$wsDe.Columns.Item(3).ColumnWidth = 17
